# Einkauf Getränke neu berechnet. Schüwo mit TabBar
#
# Adds a new row of expense data ("Kaffee und Pocorn") to the Ausgaben sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Ausgaben")

# Row 19: Kategorie / Bezeichnung / Datum / Betrag / Firmennamen / Adresse
$ws.Range("A19").Value = "Kiosk"
$ws.Range("C19").Value = "Kaffee und Pocorn "
$ws.Range("D19").Value = 45316
$ws.Range("E19").Value = 32.7
$ws.Range("F19").Value = "Nadia Wagner"
$ws.Range("G19").Value = "Wilifeld 4, 5708 Birrwil"

# Reflect the active cell / selection change recorded for this sheet
$ws.Range("C19").Select()

$wb.Save()
